$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "22.013.29"
$ws.Range("E2").Value = "  -0.41%  "
$ws.Range("D3").Value = "1.553.54"
$ws.Range("E3").Value = "  +0.17%  "
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.001"
$ws.Range("E5").Value = "  +0.17%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "288.78"
$ws.Range("E6").Value = "  +0.62%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3920"
$ws.Range("E7").Value = "  +3.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3187"
$ws.Range("E8").Value = "  -2.67%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "44.33"
$ws.Range("E9").Value = "  +2.39%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07148"
$ws.Range("E10").Value = "  -2.40%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.061"
$ws.Range("E11").Value = "  -5.94%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.002"
$ws.Range("E12").Value = "  +0.18%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.613"
$ws.Range("E13").Value = "  -3.12%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.51"
$ws.Range("E14").Value = "  -7.40%  "
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "1.554.61"
$ws.Range("E15").Value = "  +1.14%  "
$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.596"
$ws.Range("E16").Value = "  -2.44%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001094"
$ws.Range("E17").Value = "  +0.63%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06552"
$ws.Range("E18").Value = "  -0.55%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "82.65"
$ws.Range("E19").Value = "  -3.34%  "
$ws.Range("E20").Value = "  +0.09%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.151"
$ws.Range("E21").Value = "  -3.40%  "
$ws.Range("E22").Value = "  -4.72%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.09"
$ws.Range("E23").Value = "  -4.75%  "
$ws.Range("B24").Value = "Toncoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.384"
$ws.Range("E24").Value = "  +4.04%  "
$ws.Range("B25").Value = "WrappedBTC"
$ws.Range("C25").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D25").Value = "22.015.41"
$ws.Range("E25").Value = "  -0.37%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.352"
$ws.Range("E26").Value = "  -6.66%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "147.00"
$ws.Range("E27").Value = "  -1.63%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.34"
$ws.Range("E28").Value = "  -3.59%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.861"
$ws.Range("E29").Value = "  -0.82%  "
$ws.Range("D30").Value = "1.727.12"
$ws.Range("E30").Value = "  +0.56%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "117.14"
$ws.Range("E31").Value = "  -3.21%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9690"
$ws.Range("E32").Value = "  -8.85%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.825"
$ws.Range("E33").Value = "  -0.77%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08240"
$ws.Range("E34").Value = "  +0.38%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.049"
$ws.Range("E35").Value = "  -2.20%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.596"
$ws.Range("E36").Value = "  -14.35%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02233"
$ws.Range("E37").Value = "  -3.23%  "
$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05976"
$ws.Range("E38").Value = "  -3.84%  "
$ws.Range("B39").Value = "InternetComputer(DFINITY)"
$ws.Range("C39").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.050"
$ws.Range("E39").Value = "  -3.89%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.200"
$ws.Range("E40").Value = "  -3.50%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2021"
$ws.Range("E41").Value = "  -5.90%  "
$ws.Range("E42").Value = "  +0.09%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "10.60"
$ws.Range("E43").Value = "  -3.36%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5731"
$ws.Range("E44").Value = "  -4.60%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.737"
$ws.Range("E45").Value = "  +0.27%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.84"
$ws.Range("E46").Value = "  -5.44%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5500"
$ws.Range("E47").Value = "  -5.20%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "116.50"
$ws.Range("E48").Value = "  -4.34%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.858"
$ws.Range("E49").Value = "  -5.92%  "
$ws.Range("E50").Value = "  -4.07%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06797"
$ws.Range("E51").Value = "  -3.05%  "
